$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text / timestamp updates (status moves from "In Translation" to "Ready for handoff") ---

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the "Latest HO Xliff Generate Date" (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 01:05:30"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 01:05:25"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2) -- this cell shared the same
# underlying string as Overview!G2 in the original workbook ("2016-09-05 01:04:53"), so it
# moves together with it to the new value ("2016-09-05 01:05:30").
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 01:05:30"

# --- Column width updates (Status column widened to fit the new, longer text) ---
# Target widths from the new content ("Ready for handoff") auto-sizing.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.3   # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3        # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3        # column C
